$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" shared text
# changes from 2016-08-25 10:15:49 to 2016-08-25 10:16:42 (shared by Overview.G and de-de.H)
$wsOverview.Range("G2").Value = "2016-08-25 10:16:42"
$wsOverview.Range("G3").Value = "2016-08-25 10:16:42"
$wsDeDe.Range("H2").Value = "2016-08-25 10:16:42"
$wsDeDe.Range("H3").Value = "2016-08-25 10:16:42"

# "Priority" value changes from "ht" to "mt" (shared by zh-cn.E and de-de.E)
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E3").Value = "mt"
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E3").Value = "mt"

# zh-cn "Correspond Handoff Datetime" changes from 2016-08-25 10:15:44 to 2016-08-25 10:16:36
$wsZhCn.Range("H2").Value = "2016-08-25 10:16:36"
$wsZhCn.Range("H3").Value = "2016-08-25 10:16:36"

# zh-cn "Correspond Handback DateTime" changes from 2016-08-25 10:16:02 to 2016-08-25 10:16:54
$wsZhCn.Range("K2").Value = "2016-08-25 10:16:54"
$wsZhCn.Range("K3").Value = "2016-08-25 10:16:54"

# de-de "Correspond Handback DateTime" changes from 2016-08-25 10:16:14 to 2016-08-25 10:17:05
$wsDeDe.Range("K2").Value = "2016-08-25 10:17:05"
$wsDeDe.Range("K3").Value = "2016-08-25 10:17:05"
